$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates (Target cluster becomes "ECs" due to new cluster, and values recomputed)
$ws.Range("D2").Value = "ECs"
$ws.Range("H2").Value = 0.601357
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.2385553333333333
$ws.Range("N2").Value = 0.7156659999999999
$ws.Range("O2").Value = 0.9144827885830529
$ws.Range("P2").Value = 0.914482788583053
$ws.Range("Q2").Value = 0.04781897319577778
$ws.Range("R2").Value = 0.430370758762
$ws.Range("S2").Value = 0.9144827885830529
$ws.Range("T2").Value = 0.914482788583053

# New row 3: same FAPs/Tac2/Tacr2 combination, but Target cluster = MuSCs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Tac2"
$ws.Range("C3").Value = "Tacr2"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.2004523333333333
$ws.Range("H3").Value = 0.601357
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.02230833333333333
$ws.Range("N3").Value = 0.066925
$ws.Range("O3").Value = 0.08551721141694704
$ws.Range("P3").Value = 0.08551721141694704
$ws.Range("Q3").Value = 0.004471757469444445
$ws.Range("R3").Value = 0.040245817225
$ws.Range("S3").Value = 0.08551721141694704
$ws.Range("T3").Value = 0.08551721141694704
